# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement database was updated: the employee
# "YIRIS FERNANDEZ CABARCAS" (row 17) was removed from this company's
# (NIT 9001120575) sheet, leaving only "MARIA DE LOS ANGELES RODRIGUEZ
# DOMINGUEZ" (row 16) — i.e. "parte 1" of the new account statement.
# The summary totals at the top of the sheet are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second worker's data row entirely (shifts every row below it
# up by one, which also re-packs the shared-string table / row numbers,
# matching how Excel itself handles a row delete).
$ws.Rows(17).Delete()

# Update the summary figures so they reflect only the remaining worker:
# "VALOR MORA" total (was 36341 + 64940 = 101281, now just 36341)
$ws.Range("E11").Value = 36341
# "Cant. Trabajadores" (worker count): 2 -> 1
$ws.Range("C13").Value = 1
# "Cant. Periodos" (period count): 2 -> 1
$ws.Range("F13").Value = 1
